# ---------------------------------------------------------------------------
# Applies the "PO Forecast" sheet addition + header renames described by the
# diff:
#   * "Weekly Quantity" sheet: B1 header "Requested quantity" -> "Weekly_PO_Qty"
#   * "Monthly Trend" sheet:  B1 header "Requested quantity" -> "Monthly_PO_Qty"
#   * New "PO Forecast" worksheet appended at the end with columns
#     ds / PO_Forecast / yhat_lower / yhat_upper and 34 rows of forecast data.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on the existing sheets -------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the look of the "Weekly Quantity" sheet so the new sheet gets the
# same header style (bold / centered / bordered) and date number format
# that the other two sheets already use.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A35").PasteSpecial(-4122)

# --- 3. Header row -----------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4. Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper) -------
$data = @(
    ,@(44976.99999999999, 4, -8.863752217277275, 17.37657567809445)
    ,@(44990.99999999999, 5, -8.163231164931091, 18.61714131979594)
    ,@(45004.99999999999, 5, -8.407632895344351, 17.83359362483577)
    ,@(45025.99999999999, 5, -8.368502735021089, 17.34483663378481)
    ,@(45081.99999999999, 5, -8.120147473852537, 18.30896370684232)
    ,@(45144.99999999999, 6, -8.322710530751458, 18.73806387460757)
    ,@(45151.99999999999, 6, -7.531188832294291, 19.96388015454887)
    ,@(45179.99999999999, 6, -5.893636130882616, 19.93481941671612)
    ,@(45207.99999999999, 7, -6.441859827070132, 19.71117495733529)
    ,@(45221.99999999999, 7, -6.13950970084772, 19.71233738974309)
    ,@(45263.99999999999, 7, -6.246932920252646, 20.48104428389698)
    ,@(45277.99999999999, 7, -6.890831424946911, 20.42499782385294)
    ,@(45305.99999999999, 8, -5.942304693478872, 22.0765607807224)
    ,@(45361.99999999999, 8, -4.834962518153404, 21.82681247209772)
    ,@(45389.99999999999, 8, -5.978748154826236, 21.82443656366103)
    ,@(45417.99999999999, 9, -4.362368979307549, 22.06869720504935)
    ,@(45438.99999999999, 9, -4.064081976801312, 22.7242228554031)
    ,@(45459.99999999999, 9, -4.164127685522659, 22.33872902696784)
    ,@(45466.99999999999, 9, -3.99565414927458, 22.69796009253515)
    ,@(45480.99999999999, 9, -3.791389606823153, 23.39827442493442)
    ,@(45494.99999999999, 10, -4.426363668811764, 23.2421798351952)
    ,@(45501.99999999999, 10, -4.474503323858724, 22.4822476758218)
    ,@(45522.99999999999, 10, -3.702339010122818, 22.56193887567983)
    ,@(45529.99999999999, 10, -2.518460038882307, 23.78225727223554)
    ,@(45585.99999999999, 10, -3.469194994225089, 23.6124586471685)
    ,@(45592.99999999999, 10, -2.885475145194926, 23.5461290780766)
    ,@(45599.99999999999, 11, -2.572425221434706, 24.26699283920285)
    ,@(45606.99999999999, 11, -2.882497665683946, 25.6775763419128)
    ,@(45613.99999999999, 11, -2.093138833483065, 24.26546576517479)
    ,@(45620.99999999999, 11, -3.362053852023168, 23.59810979921912)
    ,@(45627.99999999999, 11, -2.88092926262949, 24.31263551686187)
    ,@(45634.99999999999, 11, -3.12170180644169, 22.93194986834898)
    ,@(45641.99999999999, 11, -2.812393629927346, 24.94809426241415)
    ,@(45648.99999999999, 11, -1.974198759088988, 24.09543082893821)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

Write-Host "PO Forecast sheet created with $($r - 2) data rows."
